$wb = $excel.ActiveWorkbook
$wsBlocks = $wb.Worksheets.Item("Workblocks")

# Remove the "wbCloseAppsRecover" workblock pair (rows 5:6) - its sibling
# workflow entries (CloseAllApplications / InitAllApplications / Process)
# no longer ship a default workblock, so the rows below shift up and the
# trailing rows are cleared instead of holding "CloseApps"/"InitApps"/"ProcessApps".
$wsBlocks.Rows("5:6").Delete()

# Clear the now-stale CloseAllApplications / InitAllApplications / Process
# workblock rows (post-shift rows 9:14), keeping the left-aligned formatting
# on column B that the template used for workblock value cells.
$wsBlocks.Range("A9:C14").ClearContents()
$wsBlocks.Range("B9:B14").HorizontalAlignment = -4131

# Reselect to match the new layout and make "Workblocks" the active sheet/tab.
$wsBlocks.Range("A3:C8").Select()
$wsBlocks.Activate()
